$wb = $excel.ActiveWorkbook

# --- 1. Update Status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: Status columns are E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column is C, row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column is C, row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the Status columns' width ---
# Target stored width in the workbook XML is 13.4101845877511 (down from
# 17.2159881591797). The engine quantizes ColumnWidth to an MDW-7 pixel
# grid, so we pick the input value whose quantized result lands as close
# as possible to the target (13.333333333333334).
$newWidth = 12.55

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
